# Add 20 new "quiz" rows (content/answer pairs) to Sheet1, rows 34-53,
# and move the view/selection to the bottom of the new data - mirroring
# the upstream "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows -----------------------------------------------------
# Column A = question text (becomes a new shared string), column B =
# numeric answer. Rows are written in ascending row order EXCEPT 50/51:
# the source workbook's shared-string table has index 50 ("so lon nhat
# co mot chu so") ahead of index 51 ("so lien sau cua so tam") even
# though row 50 references string 51 and row 51 references string 50.
# Writing A51 before A50 reproduces that exact shared-string ordering.

$ws.Range("A34").Value = "số lẻ liền sau số ba"
$ws.Range("B34").Value = 5

$ws.Range("A35").Value = "số nguyên tố nhỏ thứ 3"
$ws.Range("B35").Value = 5

$ws.Range("A36").Value = "một nửa của một chục"
$ws.Range("B36").Value = 5

$ws.Range("A37").Value = "số liền sau của số năm"
$ws.Range("B37").Value = 6

$ws.Range("A38").Value = "số chẵn lớn thứ hai có một chữ số"
$ws.Range("B38").Value = 6

$ws.Range("A39").Value = "Trong một năm có bao nhiêu tháng có 30 ngày"
$ws.Range("B39").Value = 6

$ws.Range("A40").Value = "số liền trước của số 7"
$ws.Range("B40").Value = 6

$ws.Range("A41").Value = "một tuần có mấy ngày"
$ws.Range("B41").Value = 7

$ws.Range("A42").Value = "số lẻ lớn thứ hai có một chữ số"
$ws.Range("B42").Value = 7

$ws.Range("A43").Value = "số liến trước của số tám"
$ws.Range("B43").Value = 7

$ws.Range("A44").Value = "số liền sau của số sáu"
$ws.Range("B44").Value = 7

$ws.Range("A45").Value = "số chẵn lớn nhất có một chữ số"
$ws.Range("B45").Value = 8

$ws.Range("A46").Value = "số liến trước của số chín"
$ws.Range("B46").Value = 8

$ws.Range("A47").Value = "số liền sau của số bảy"
$ws.Range("B47").Value = 8

$ws.Range("A48").Value = "tháng chẵn đầu tiên có 31 ngày"
$ws.Range("B48").Value = 8

$ws.Range("A49").Value = "số lẻ lớn nhất có một chữ số"
$ws.Range("B49").Value = 9

# Row 51 written first on purpose (see note above).
$ws.Range("A51").Value = "số lớn nhất có một chữ số"
$ws.Range("B51").Value = 9

$ws.Range("A50").Value = "số liền sau của số tám"
$ws.Range("B50").Value = 9

$ws.Range("A52").Value = "số liền trước của số mười"
$ws.Range("B52").Value = 9

$ws.Range("A53").Value = "số chính phương lớn nhất có một chữ số"
$ws.Range("B53").Value = 9

# --- View state ----------------------------------------------------------
# Scroll the window so row 24 is near the top, and leave the selection on
# the last cell that was filled in, matching the author's final view.
$win = $excel.ActiveWindow
$win.ScrollRow = 24
$win.ScrollColumn = 1
$ws.Range("B53").Select()
